$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room: insert 3 fresh blank rows just above the old SUM row (15),
#    pushing it down to row 18. After this, rows 4-14 still hold the old
#    data (untouched) and rows 15-17 are blank.
# ---------------------------------------------------------------------------
$ws.Rows("15:17").Insert()

# ---------------------------------------------------------------------------
# 2. Wipe every data cell B4:H18 so we can re-author the whole BOM table
#    from scratch without worrying about stale leftovers from the old
#    row layout.
# ---------------------------------------------------------------------------
$ws.Range("B4:H18").ClearContents()

# ---------------------------------------------------------------------------
# 3. Re-author the table content, row by row, exactly as it should read
#    after the edit.
# ---------------------------------------------------------------------------

# Row 4 - ATMEGA32u4-AU (unchanged)
$ws.Range("B4").Value = "mouser"
$ws.Range("C4").Value = "http://www.mouser.ch/ProductDetail/Atmel/ATMEGA32U4-AU/?qs=sGAEpiMZZMvAvBNgSS9Lqh1wBIpnYT9c"
$ws.Range("D4").Value = "ATMEGA32u4-AU"
$ws.Range("E4").Value = "8bits 32k cpu"
$ws.Range("F4").Value = 6.33
$ws.Range("G4").Value = 1
$ws.Range("H4").Formula = "=F4*G4"

# Row 5 - CSTCE16M0V53ZW-R0 resonator (same data, total formula now present)
$ws.Range("B5").Value = "mouser"
$ws.Range("C5").Value = "http://www.mouser.ch/Search/ProductDetail.aspx?qs=8%2f1pEl6ptNseo9Gxrhu%2fPA%3d%3d"
$ws.Range("D5").Value = "CSTCE16M0V53ZW-R0"
$ws.Range("E5").Value = "resonator 16Mhz"
$ws.Range("F5").Value = 0.5
$ws.Range("G5").Value = 1
$ws.Range("H5").Formula = "=F5*G5"
$ws.Rows("5:5").RowHeight = 15.65

# Row 6 - MCP73831 LiPo charger (unchanged)
$ws.Range("B6").Value = "mouser"
$ws.Range("C6").Value = "http://www.mouser.ch/ProductDetail/Microchip-Technology/MCP73831T-2ACI-OT/?qs=sGAEpiMZZMtLck3p7ZBovc%252bIEf4wKPGR"
$ws.Range("D6").Value = "MCP73831"
$ws.Range("E6").Value = "MCP73831 LiPo charger"
$ws.Range("F6").Value = 0.6
$ws.Range("G6").Value = 1
$ws.Range("H6").Formula = "=F6*G6"

# Row 7 - NEW: MAX1595EUA33 buck/boost 3.3v dcdc
$ws.Range("B7").Value = "mouser"
$ws.Range("D7").Value = "MAX1595EUA33"
$ws.Range("E7").Value = "buck/boost 3.3v dcdc"
$ws.Range("F7").Value = 3.88
$ws.Range("G7").Value = 1
$ws.Range("H7").Formula = "=F7*G7"
$ws.Hyperlinks.Add($ws.Range("C7"), "http://www.mouser.ch/ProductDetail/Maxim-Integrated/MAX1595EUA33+/?qs=sGAEpiMZZMtitjHzVIkrqUmW7fHvDhXHgnQoEKfsHaU%3d", "", "", "http://www.mouser.ch/ProductDetail/Maxim-Integrated/MAX1595EUA33+/?qs=sGAEpiMZZMtitjHzVIkrqUmW7fHvDhXHgnQoEKfsHaU%3d")
$ws.Range("D6").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Rows("7:7").RowHeight = 15.65

# Row 8 - conn USB Micro B femal (unchanged, shifted down)
$ws.Range("B8").Value = "mouser"
$ws.Range("E8").Value = "conn USB Micro B femal"
$ws.Range("F8").Value = 0.7
$ws.Range("G8").Value = 1
$ws.Range("H8").Formula = "=F8*G8"

# Row 9 - SMT-1141-T-3-R piezzo (unchanged, shifted down)
$ws.Range("B9").Value = "mouser"
$ws.Range("C9").Value = "http://www.mouser.ch/ProductDetail/PUI-Audio/SMT-1141-T-3-R/?qs=%2fha2pyFaduiCLgby5iJPqp2iGH%252b6CT48ZBFKPYK%2fGYsaasTvnsIatg%3d%3d"
$ws.Range("D9").Value = "SMT-1141-T-3-R"
$ws.Range("E9").Value = "piezzo "
$ws.Range("F9").Value = 3.9
$ws.Range("G9").Value = 1
$ws.Range("H9").Formula = "=F9*G9"

# Row 10 - SKQGAKE010 push button (price per-unit corrected)
$ws.Range("B10").Value = "mouser"
$ws.Range("C10").Value = "http://www.mouser.ch/ProductDetail/ALPS/SKQGAKE010/?qs=sGAEpiMZZMtFyPk3yBMYYCtu4vPfeeUaHHDNk5wDwBc%3d"
$ws.Range("D10").Value = "SKQGAKE010"
$ws.Range("E10").Value = "push button"
$ws.Range("F10").Value = 0.16
$ws.Range("G10").Value = 6
$ws.Range("H10").Formula = "=F10*G10"

# Row 11 - NEW: SK-12C0405-SG 1.5 RT switch on-off
$ws.Range("B11").Value = "mouser"
$ws.Range("D11").Value = "SK-12C0405-SG 1.5 RT"
$ws.Range("E11").Value = "switch on-off"
$ws.Range("F11").Value = 0.97
$ws.Range("G11").Value = 1
$ws.Range("H11").Formula = "=F11*G11"
$ws.Hyperlinks.Add($ws.Range("C11"), "http://www.mouser.ch/ProductDetail/CK-Components/SK-12C0405-SG-15-RT/?qs=sGAEpiMZZMtHXLepoqNyVaknRufv4Zo6J8yLuspm3Zw%3d", "", "", "http://www.mouser.ch/ProductDetail/CK-Components/SK-12C0405-SG-15-RT/?qs=sGAEpiMZZMtHXLepoqNyVaknRufv4Zo6J8yLuspm3Zw%3d")
$ws.Range("D6").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Rows("11:11").RowHeight = 15.65

# Row 12 - NEW: KMR631NG ULC LFS reset push buton
$ws.Range("B12").Value = "mouser"
$ws.Range("D12").Value = "KMR631NG ULC LFS"
$ws.Range("E12").Value = "reset push buton"
$ws.Range("F12").Value = 0.27
$ws.Range("G12").Value = 1
$ws.Range("H12").Formula = "=F12*G12"
$ws.Hyperlinks.Add($ws.Range("C12"), "http://www.mouser.ch/ProductDetail/CK-Components/KMR631NG-ULC-LFS/?qs=sGAEpiMZZMsgGjVA3toVBJ1OkFFtNMGB4KijNZUSro0%3d", "", "", "http://www.mouser.ch/ProductDetail/CK-Components/KMR631NG-ULC-LFS/?qs=sGAEpiMZZMsgGjVA3toVBJ1OkFFtNMGB4KijNZUSro0%3d")
$ws.Range("D6").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Rows("12:12").RowHeight = 15.7

# Row 13 - led RGB (unchanged, shifted down)
$ws.Range("B13").Value = "mouser"
$ws.Range("E13").Value = "led RGB"
$ws.Range("F13").Value = 2.5
$ws.Range("G13").Value = 1
$ws.Range("H13").Formula = "=F13*G13"

# Row 14 - NEW: BAT-HLD-002-SMT support CR2016
$ws.Range("B14").Value = "mouser"
$ws.Range("D14").Value = "BAT-HLD-002-SMT"
$ws.Range("E14").Value = "support CR2016"
$ws.Range("F14").Value = 0.28
$ws.Range("G14").Value = 1
$ws.Range("H14").Formula = "=F14*G14"
$ws.Hyperlinks.Add($ws.Range("C14"), "http://www.mouser.ch/ProductDetail/Linx-Technologies/BAT-HLD-002-SMT/?qs=%2fha2pyFaduilhNkyJFgy2WekJWQQ7JGY1Lox0Z3adM0%3d", "", "", "http://www.mouser.ch/ProductDetail/Linx-Technologies/BAT-HLD-002-SMT/?qs=%2fha2pyFaduilhNkyJFgy2WekJWQQ7JGY1Lox0Z3adM0%3d")
$ws.Range("D6").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Rows("14:14").RowHeight = 15.65

# Row 15 - GEB014461 LiPo 3.6V 180 mAh (unchanged, shifted down)
$ws.Range("B15").Value = "alibaba"
$ws.Range("C15").Value = "https://www.alibaba.com/product-detail/2016-hot-sale-lithium-polymer-battery_60437032979.html"
$ws.Range("D15").Value = "GEB014461"
$ws.Range("E15").Value = "LiPo 3.6V 180 mAh"
$ws.Range("F15").Value = 5
$ws.Range("G15").Value = 1
$ws.Range("H15").Formula = "=F15*G15"

# Row 16 - ER-OLED0.96-1 (unchanged, shifted down)
$ws.Range("B16").Value = "buydisplay"
$ws.Range("C16").Value = "http://www.buydisplay.com/default/datasheet-128x64-oled-module-spi-0-96-inch-graphic-displays-white-on-black"
$ws.Range("D16").Value = "ER-OLED0.96-1"
$ws.Range("E16").Value = "oled white on black 128x64 0.96''"
$ws.Range("F16").Value = 3.61
$ws.Range("G16").Value = 1
$ws.Range("H16").Formula = "=F16*G16"

# Row 17 - seeed PCB (unchanged, shifted down)
$ws.Range("B17").Value = "seeed"
$ws.Range("D17").Value = "PCB"
$ws.Range("E17").Value = "pcb 85.6x53.98 x 0.6mm 2 layers black"
$ws.Range("F17").Value = 2
$ws.Range("G17").Value = 1
$ws.Range("H17").Formula = "=F17*G17"

# Row 18 - grand total
$ws.Range("H18").Formula = "=SUM(H4:H17)"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Cosmetic sheet-view touch ups (column widths & selection) to match the
#    edited workbook as closely as possible.
# ---------------------------------------------------------------------------
$ws.Columns("A:B").ColumnWidth = 10.7
$ws.Columns("C:C").ColumnWidth = 10.4
$ws.Columns("D:D").ColumnWidth = 20.8
$ws.Columns("E:E").ColumnWidth = 40.2
$ws.Columns("F:F").ColumnWidth = 10.7

$ws.Range("H22").Select()
